$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column E
$ws.Range("E1").Value = "aware of sound presentation?"

# Match the header formatting used by the other header cells (bold, fill, centered)
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Widen column E to fit the new header text
$ws.Range("E:E").ColumnWidth = 26.33203125

# Data values for E2:E73 (1 = yes, 0 = no)
$values = @(0,0,0,0,0,0,0,0,0,0,0,1,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,1,0,0,1,0,0,1,0,0,0,0,0,1,0,1,1,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 5).Value = $values[$i]
}

# Summary formula in E75
$ws.Range("E75").Formula = "=(SUM(E2:E74)/72)*100"

# Threaded comment on E1
$comment = $ws.Range("E1").AddCommentThreaded("1= yes; 0= no")

$excel.DisplayAlerts = $false
